# Apply the updated cryptocurrency price/volume figures from the latest
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.343.69"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").Value = "1.870.69"
$ws.Range("E3").Value = "  +0.66%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.66"
$ws.Range("E5").Value = "  +1.48%  "

# Row 6
$ws.Range("E6").Value = "  -0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.51%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2847"
$ws.Range("E8").Value = "  +1.54%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06547"
$ws.Range("E9").Value = "  +0.34%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.50"
$ws.Range("E10").Value = "  +7.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07886"
$ws.Range("E11").Value = "  +1.88%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.15"
$ws.Range("E12").Value = "  +2.42%  "

# Row 13
$ws.Range("D13").Value = "1.871.10"
$ws.Range("E13").Value = "  +0.60%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.114"
$ws.Range("E14").Value = "  +1.49%  "

# Row 15
$ws.Range("E15").Value = "  +1.79%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.46"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("D17").Value = "30.347.13"
$ws.Range("E17").Value = "  +0.95%  "

# Row 19
$ws.Range("E19").Value = "  +2.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.475"
$ws.Range("E20").Value = "  +3.11%  "

# Row 21
$ws.Range("D21").Value = "2.115.39"
$ws.Range("E21").Value = "  -0.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007320"
$ws.Range("E22").Value = "  +1.80%  "

# Row 23
$ws.Range("E23").Value = "  -0.18%  "

# Row 24
$ws.Range("E24").Value = "  +0.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.33"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.175"
$ws.Range("E26").Value = "  -0.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.23"
$ws.Range("E27").Value = "  +2.04%  "

# Row 28
$ws.Range("E28").Value = "  +0.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.382"
$ws.Range("E29").Value = "  +1.68%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09704"
$ws.Range("E30").Value = "  +0.88%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.406"
$ws.Range("E31").Value = "  +1.60%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.475"
$ws.Range("E32").Value = "  +1.00%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.102"
$ws.Range("E33").Value = "  +0.97%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04715"
$ws.Range("E34").Value = "  +1.83%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  +4.92%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7075"
$ws.Range("E36").Value = "  +1.89%  "

# Row 37
$ws.Range("E37").Value = "  +0.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01863"
$ws.Range("E38").Value = "  +0.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.340"
$ws.Range("E39").Value = "  +0.11%  "

# Row 40
$ws.Range("E40").Value = "  +1.49%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.65"
$ws.Range("E41").Value = "  +5.35%  "

# Row 42
$ws.Range("E42").Value = "  +1.44%  "

# Row 43
$ws.Range("E43").Value = "  -0.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4194"
$ws.Range("E44").Value = "  +1.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.78"
$ws.Range("E46").Value = "  +1.11%  "

# Row 47
$ws.Range("E47").Value = "  +1.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.297"
$ws.Range("E48").Value = "  +2.82%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "940.34"
$ws.Range("E49").Value = "  -3.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.29"
$ws.Range("E50").Value = "  +1.94%  "

# Row 51
$ws.Range("E51").Value = "  -0.58%  "
